$wb = $excel.ActiveWorkbook

# --- Summary sheet: update NPV, add new surplus/unmet-demand rows ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 2896166.6806224
$wsSummary.Range("A7").Value = "Wasted Prosumer Surplus"
$wsSummary.Range("B7").Value = 342857.5917216506
$wsSummary.Range("A8").Value = "Total Wasted Prosumer Surplus"
$wsSummary.Range("B8").Value = 342857.5917216506
$wsSummary.Range("A9").Value = "Unmet Demand"
$wsSummary.Range("B9").Value = 29348.32150215295
$wsSummary.Range("A10").Value = "Total Unmet Demand"
$wsSummary.Range("B10").Value = 29348.32150215295
$wsSummary.Range("A6").Copy()
$wsSummary.Range("A7:A10").PasteSpecial(-4122)

# --- Costs and Revenues sheet ---
$ws_Costs_and_Revenues = $wb.Worksheets.Item("Costs and Revenues")
$ws_Costs_and_Revenues.Range("E2").Value = 938737.9016593838
$ws_Costs_and_Revenues.Range("F2").Value = 958748.0428165476
$ws_Costs_and_Revenues.Range("G2").Value = 959349.9008549106
$ws_Costs_and_Revenues.Range("I2").Value = 959349.9008549106
$ws_Costs_and_Revenues.Range("J2").Value = 934889.8449876352
$ws_Costs_and_Revenues.Range("K2").Value = 961571.7997246225
$ws_Costs_and_Revenues.Range("L2").Value = 969896.8330247513
$ws_Costs_and_Revenues.Range("N2").Value = 969896.8330247513
$ws_Costs_and_Revenues.Range("O2").Value = 916492.3349876349
$ws_Costs_and_Revenues.Range("P2").Value = 840347.3210590899
$ws_Costs_and_Revenues.Range("K3").Value = 47200
$ws_Costs_and_Revenues.Range("L3").Value = 19200
$ws_Costs_and_Revenues.Range("B4").Value = 516729.2344596348
$ws_Costs_and_Revenues.Range("C4").Value = 514932.0834136077
$ws_Costs_and_Revenues.Range("D4").Value = 513132.4944200165
$ws_Costs_and_Revenues.Range("E4").Value = 501136.9470735459
$ws_Costs_and_Revenues.Range("F4").Value = 511568.7470510048
$ws_Costs_and_Revenues.Range("G4").Value = 510118.6542394872
$ws_Costs_and_Revenues.Range("H4").Value = 508299.2497951454
$ws_Costs_and_Revenues.Range("I4").Value = 506477.3054326685
$ws_Costs_and_Revenues.Range("J4").Value = 489991.0232228274
$ws_Costs_and_Revenues.Range("K4").Value = 504151.9865617157
$ws_Costs_and_Revenues.Range("L4").Value = 507265.0393204461
$ws_Costs_and_Revenues.Range("M4").Value = 505406.1812562778
$ws_Costs_and_Revenues.Range("N4").Value = 503544.650463314
$ws_Costs_and_Revenues.Range("O4").Value = 470341.4757097192
$ws_Costs_and_Revenues.Range("P4").Value = 424118.69410044
$ws_Costs_and_Revenues.Range("K5").Value = 34888.635
$ws_Costs_and_Revenues.Range("P5").Value = 25472.907
$ws_Costs_and_Revenues.Range("B6").Value = 400618.8526664367
$ws_Costs_and_Revenues.Range("C6").Value = 406416.0037124637
$ws_Costs_and_Revenues.Range("D6").Value = 408215.5927060549
$ws_Costs_and_Revenues.Range("E6").Value = 101234.3895858378
$ws_Costs_and_Revenues.Range("F6").Value = 392711.0057655428
$ws_Costs_and_Revenues.Range("G6").Value = 409878.8876154234
$ws_Costs_and_Revenues.Range("H6").Value = 416498.2920597651
$ws_Costs_and_Revenues.Range("I6").Value = 418320.2364222421
$ws_Costs_and_Revenues.Range("J6").Value = 132868.5327648078
$ws_Costs_and_Revenues.Range("K6").Value = 375331.1781629068
$ws_Costs_and_Revenues.Range("L6").Value = 407029.9167043052
$ws_Costs_and_Revenues.Range("M6").Value = 428088.7747684736
$ws_Costs_and_Revenues.Range("N6").Value = 429950.3055614373
$ws_Costs_and_Revenues.Range("O6").Value = 192517.8812779157
$ws_Costs_and_Revenues.Range("P6").Value = 390755.7199586499

# --- Installed Capacities sheet ---
$ws_Installed_Capacities = $wb.Worksheets.Item("Installed Capacities")
$ws_Installed_Capacities.Range("K2").Value = 415
$ws_Installed_Capacities.Range("P2").Value = 303

# --- Added Capacities sheet ---
$ws_Added_Capacities = $wb.Worksheets.Item("Added Capacities")
$ws_Added_Capacities.Range("K2").Value = 59
$ws_Added_Capacities.Range("L2").Value = 24

# --- Retired Capacities sheet ---
$ws_Retired_Capacities = $wb.Worksheets.Item("Retired Capacities")
$ws_Retired_Capacities.Range("P2").Value = 59

# --- DG Dispatch sheet ---
$ws_DG_Dispatch = $wb.Worksheets.Item("DG Dispatch")
$ws_DG_Dispatch.Range("B29").Value = 415
$ws_DG_Dispatch.Range("C29").Value = 415
$ws_DG_Dispatch.Range("H29").Value = 415
$ws_DG_Dispatch.Range("T29").Value = 415
$ws_DG_Dispatch.Range("U29").Value = 415
$ws_DG_Dispatch.Range("V29").Value = 415
$ws_DG_Dispatch.Range("W29").Value = 415
$ws_DG_Dispatch.Range("X29").Value = 415
$ws_DG_Dispatch.Range("Y29").Value = 415
$ws_DG_Dispatch.Range("W30").Value = 415
$ws_DG_Dispatch.Range("X30").Value = 415
$ws_DG_Dispatch.Range("Q31").Value = 415
$ws_DG_Dispatch.Range("R31").Value = 415
$ws_DG_Dispatch.Range("B44").Value = 303
$ws_DG_Dispatch.Range("C44").Value = 303
$ws_DG_Dispatch.Range("D44").Value = 303
$ws_DG_Dispatch.Range("E44").Value = 303
$ws_DG_Dispatch.Range("F44").Value = 303
$ws_DG_Dispatch.Range("G44").Value = 303
$ws_DG_Dispatch.Range("H44").Value = 303
$ws_DG_Dispatch.Range("T44").Value = 303
$ws_DG_Dispatch.Range("U44").Value = 303
$ws_DG_Dispatch.Range("V44").Value = 303
$ws_DG_Dispatch.Range("W44").Value = 303
$ws_DG_Dispatch.Range("X44").Value = 303
$ws_DG_Dispatch.Range("Y44").Value = 303
$ws_DG_Dispatch.Range("B45").Value = 303
$ws_DG_Dispatch.Range("C45").Value = 303
$ws_DG_Dispatch.Range("D45").Value = 303
$ws_DG_Dispatch.Range("E45").Value = 303
$ws_DG_Dispatch.Range("F45").Value = 303
$ws_DG_Dispatch.Range("G45").Value = 303
$ws_DG_Dispatch.Range("H45").Value = 303
$ws_DG_Dispatch.Range("R45").Value = 303
$ws_DG_Dispatch.Range("S45").Value = 303
$ws_DG_Dispatch.Range("T45").Value = 303
$ws_DG_Dispatch.Range("U45").Value = 303
$ws_DG_Dispatch.Range("V45").Value = 303
$ws_DG_Dispatch.Range("W45").Value = 303
$ws_DG_Dispatch.Range("X45").Value = 303
$ws_DG_Dispatch.Range("Y45").Value = 303
$ws_DG_Dispatch.Range("Q46").Value = 303
$ws_DG_Dispatch.Range("R46").Value = 303
$ws_DG_Dispatch.Range("S46").Value = 303

# --- Unmet Demand sheet ---
$ws_Unmet_Demand = $wb.Worksheets.Item("Unmet Demand")
$ws_Unmet_Demand.Range("B29").Value = 66.99931295557451
$ws_Unmet_Demand.Range("C29").Value = 34.47457824299391
$ws_Unmet_Demand.Range("H29").Value = 20.02773927029563
$ws_Unmet_Demand.Range("T29").Value = 114.6191915811053
$ws_Unmet_Demand.Range("U29").Value = 229.5217529288726
$ws_Unmet_Demand.Range("V29").Value = 214.8510241668239
$ws_Unmet_Demand.Range("W29").Value = 223.3734759809475
$ws_Unmet_Demand.Range("X29").Value = 177.2818334606677
$ws_Unmet_Demand.Range("Y29").Value = 96.31743268280638
$ws_Unmet_Demand.Range("W30").Value = 17.37314290982852
$ws_Unmet_Demand.Range("X30").Value = 4.862739445387547
$ws_Unmet_Demand.Range("Q31").Value = 107.1821235684552
$ws_Unmet_Demand.Range("R31").Value = 306.1956210454637
$ws_Unmet_Demand.Range("B44").Value = 178.9993129555745
$ws_Unmet_Demand.Range("C44").Value = 146.4745782429939
$ws_Unmet_Demand.Range("D44").Value = 107.3391557398498
$ws_Unmet_Demand.Range("E44").Value = 101.3632896068686
$ws_Unmet_Demand.Range("F44").Value = 101.8896287080119
$ws_Unmet_Demand.Range("G44").Value = 107.8327491714383
$ws_Unmet_Demand.Range("H44").Value = 132.0277392702956
$ws_Unmet_Demand.Range("T44").Value = 226.6191915811053
$ws_Unmet_Demand.Range("U44").Value = 341.5217529288726
$ws_Unmet_Demand.Range("V44").Value = 326.8510241668239
$ws_Unmet_Demand.Range("W44").Value = 335.3734759809475
$ws_Unmet_Demand.Range("X44").Value = 289.2818334606677
$ws_Unmet_Demand.Range("Y44").Value = 208.3174326828064
$ws_Unmet_Demand.Range("B45").Value = 81.55655664632661
$ws_Unmet_Demand.Range("C45").Value = 58.09991244551929
$ws_Unmet_Demand.Range("D45").Value = 44.93768689770263
$ws_Unmet_Demand.Range("E45").Value = 39.67209722191262
$ws_Unmet_Demand.Range("F45").Value = 36.63624233787687
$ws_Unmet_Demand.Range("G45").Value = 26.52519625238585
$ws_Unmet_Demand.Range("H45").Value = 42.22842014979517
$ws_Unmet_Demand.Range("R45").Value = 50.67054165050009
$ws_Unmet_Demand.Range("S45").Value = 37.140588939824
$ws_Unmet_Demand.Range("T45").Value = 85.53401876295709
$ws_Unmet_Demand.Range("U45").Value = 96.68869740971195
$ws_Unmet_Demand.Range("V45").Value = 111.5106671915202
$ws_Unmet_Demand.Range("W45").Value = 129.3731429098285
$ws_Unmet_Demand.Range("X45").Value = 116.8627394453875
$ws_Unmet_Demand.Range("Y45").Value = 96.39139276613435
$ws_Unmet_Demand.Range("Q46").Value = 219.1821235684552
$ws_Unmet_Demand.Range("R46").Value = 418.1956210454637
$ws_Unmet_Demand.Range("S46").Value = 106.541226054864
